$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 3 new rows right before the old last data row (row 8), so that:
#    - old row 8 (thick-bottom bordered data row) shifts down to row 11
#    - old row 9 (Average / thick-bottom summary row) shifts down to row 12
#    Formulas and the merged-cell range auto-adjust when rows are inserted inside them.
$ws.Rows.Item(8).Insert()
$ws.Rows.Item(8).Insert()
$ws.Rows.Item(8).Insert()

# 2. The 3 freshly-inserted rows (8,9,10) are blank. After the insert, row 11 holds the
#    formatting that used to belong to the old last data row (the "box-closing" border look).
#    Copy that formatting down onto the 3 new rows too (only row 11 keeps the thick-bottom
#    row height, since that is a row-level property that PasteSpecial(formats) won't touch).
$ws.Range("A11:L11").Copy()
$ws.Range("A8:L10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3. Overwrite the calibration values for the 10 data rows (2-11) with the new readings.
$data = @(
  @(4,   -48, -42, -20, 159, -101, -1, 1,  0, 1000, 818),
  @(-13, -47, -14, -19, 177, -114,  1, 0,  0, 1000, 994),
  @(-6,  -34, -24, -15, 156, -106, -2, 0,  0, 1000, 869),
  @(-8,  -43, -31, -34, 151, -104,  1, -1, 0, 1000, 805),
  @(6,   -38, -33, -37, 156, -108,  2, 1,  0, 1000, 818),
  @(10,  -47, -34, -33, 156, -103, -2, -2, 0, 1000, 817),
  @(-11, -43, -32, -31, 155, -108,  3, -1, 0, 1000, 798),
  @(-5,  -37, -24, -34, 163, -105, -1, 0,  0, 1000, 827),
  @(11,  -43, -14, -40, 175, -105, -2, 0,  0, 1000, 903),
  @(-7,  -33, -21, -21, 172, -102, -1, 0,  1, 1000, 881)
)

$row = 2
foreach ($vals in $data) {
    $ws.Cells.Item($row, 2).Value  = $vals[0]
    $ws.Cells.Item($row, 3).Value  = $vals[1]
    $ws.Cells.Item($row, 4).Value  = $vals[2]
    $ws.Cells.Item($row, 5).Value  = $vals[3]
    $ws.Cells.Item($row, 6).Value  = $vals[4]
    $ws.Cells.Item($row, 7).Value  = $vals[5]
    $ws.Cells.Item($row, 8).Value  = $vals[6]
    $ws.Cells.Item($row, 9).Value  = $vals[7]
    $ws.Cells.Item($row, 10).Value = $vals[8]
    $ws.Cells.Item($row, 11).Value = $vals[9]
    $ws.Cells.Item($row, 12).Value = $vals[10]
    $row = $row + 1
}

# 4. Update the selected cell shown when the workbook is reopened.
$ws.Range("N14").Select()

$wb.Application.Calculate()
